# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text (cell A1, shared string)
$ws.Range("A1").Value = "Datos actualizados a 6 de Septiembre de 2020 a las 08:11"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6431160
$ws.Range("C4").Value = 8
$ws.Range("E4").Value = 2531340
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 192820

# Row 6 - India
$ws.Range("B6").Value = 4114773
$ws.Range("C6").Value = 3934
$ws.Range("D6").Value = 3180999
$ws.Range("E6").Value = 863070
$ws.Range("G6").Value = 25
$ws.Range("H6").Value = 70704

# Row 20 - Pakistan
$ws.Range("B20").Value = 298509
$ws.Range("C20").Value = 484
$ws.Range("D20").Value = 285898
$ws.Range("E20").Value = 6269
$ws.Range("G20").Value = 2
$ws.Range("H20").Value = 6342

# Row 27 - Ucrania
$ws.Range("B27").Value = 135894
$ws.Range("C27").Value = 2107
$ws.Range("E27").Value = 71434

# Row 29 - Israel
$ws.Range("B29").Value = 129173
$ws.Range("C29").Value = 237
$ws.Range("D29").Value = 102104
$ws.Range("E29").Value = 26062

# Row 62 - Kirguistan
$ws.Range("B62").Value = 44403
$ws.Range("C62").Value = 87
$ws.Range("D62").Value = 39826
$ws.Range("E62").Value = 3517

# Row 152 - Georgia
$ws.Range("B152").Value = 1650
$ws.Range("C152").Value = 29
$ws.Range("D152").Value = 1310
$ws.Range("E152").Value = 321
